$wb = $excel.ActiveWorkbook

# --- "500 bar" sheet: Heat capacity is now a plain entered value instead of a formula ---
$ws1 = $wb.Worksheets.Item("500 bar")
$ws1.Range("B2").Value = 0.00394

# --- "LOHC_load" sheet: Heat demand (kWh per kg H2) updated, and selection moved ---
$ws3 = $wb.Worksheets.Item("LOHC_load")
$ws3.Range("B3").Value = -9
$ws3.Range("E5").Select()

# --- Re-activate the "500 bar" sheet and update its selection so it becomes the visible tab ---
$ws1.Activate()
$ws1.Range("D3").Select()
